# Updated cryptos list on Mon Feb 19 09:51:33 UTC 2024 with GitHub Actions
# Refreshes prices / 1h-volume% for each coin row, and re-sorts a few rows
# whose rank changed (Uniswap/ImmutableX, Dai/Kaspa, RenderToken/Toncoin/
# Filecoin, Stacks/Celestia/ARBITRUM).
#
# Price values that look like plain decimals (e.g. "0.633", "40.12") are
# written with a leading apostrophe so Excel stores them as text instead
# of silently re-parsing them as numbers (which would both change the
# type and lose formatting such as trailing zeros). Prices that already
# contain two dots (e.g. "52.360.39") aren't valid numbers to Excel, so
# they stay text without needing the apostrophe trick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.360.39"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.915.75"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'352.23"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'112.55"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "'40.12"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'19.97"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'7.82"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "3.375.61"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("E16").Value = "  +6.63%  "
$ws.Range("D17").Value = "2.932.09"
$ws.Range("E17").Value = "  +4.24%  "
$ws.Range("D18").Value = "52.382.67"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").Value = "'3.33"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'71.09"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'271.05"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.168"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'10.66"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("D30").Value = "'37.73"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'6.57"
$ws.Range("E31").Value = "  +5.61%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  +10.86%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.26"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  +10.61%  "
$ws.Range("D35").Value = "'53.26"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +5.86%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'18.81"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'2.07"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  +14.10%  "
$ws.Range("D42").Value = "'23.74"
$ws.Range("E42").Value = "  +7.39%  "
$ws.Range("D44").Value = "'121.56"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "'2.63"
$ws.Range("E45").Value = "  +7.27%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "'3.56"
$ws.Range("E47").Value = "  +4.62%  "
$ws.Range("D48").Value = "2.199.56"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").Value = "'0.268"
$ws.Range("E49").Value = "  +24.95%  "
$ws.Range("D50").Value = "'0.0338"
$ws.Range("E50").Value = "  +11.45%  "
$ws.Range("D51").Value = "'0.967"
$ws.Range("E51").Value = "  +3.56%  "
